$d = $word.ActiveDocument

# --- 1. Date paragraph: "16/05/20" -> three runs "21" + "/05/20" + "20",
#        with the _GoBack bookmark moved onto this paragraph. ---
$dateXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="697E5315" w14:textId="7DFE4EA0" w:rsidR="00FC28DE" w:rsidRDefault="00FC28DE" w:rsidP="00FC28DE"><w:r><w:t>21</w:t></w:r><w:r><w:t>/05/20</w:t></w:r><w:r><w:t>20</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$datePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "16/05/20`r") {
        $datePara = $p
        break
    }
}
$datePara.Range.InsertXML($dateXml)

# --- 2. Observations paragraph: rewrite the run content, drop the old
#        _GoBack bookmark (it now lives on the date paragraph). ---
$obsXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="628AC1AB" w14:textId="325CFE0D" w:rsidR="00FC28DE" w:rsidRPr="00FC28DE" w:rsidRDefault="00FC28DE" w:rsidP="00FC28DE"><w:pPr><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:color w:val="000000"/><w:szCs w:val="24"/><w:lang w:eastAsia="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Ao realizar o experimento com 80 indiv&#237;duos, utilizando a imagem </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:szCs w:val="24"/></w:rPr><w:t>16_L_N_L_BH</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:szCs w:val="24"/></w:rPr><w:t>, com 60% das imagens do indiv&#237;duo para treinamento, com o classificador SVM sem otimiza&#231;&#227;o, feature extraction utilizando SURF e k-means clustering, por meio da fun&#231;&#227;o bagOfFeatures</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:color w:val="000000"/><w:szCs w:val="24"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>, foi observada a acur&#225;cia de 98,13%.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$obsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("- Ao realizar a predição")) {
        $obsPara = $p
        break
    }
}
$obsPara.Range.InsertXML($obsXml)

# --- 3. Drop the trailing empty paragraph before the sectPr. ---
$n = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($n)
$prev = $d.Paragraphs.Item($n - 1)
$killRange = $d.Range($prev.Range.End - 1, $last.Range.End)
$killRange.Delete()
